# 自动更新Excel文件 - 每日刷新"剩余"天数（列E）
# For each data row (2..99): decrement the "剩余" (remaining days, column E)
# counter by one. When a row's counter rolls over to 0, the cycle restarts:
# "剩余" is reset to the row's "总天" (total days, column D) and the
# "开始时间" (start date, column F) is advanced by that same number of days.
#
# Row 36 already carries a fully "topped up" / stale-date record and is left
# untouched, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 99; $r++) {
    if ($r -eq 36) {
        continue
    }

    $totalDays     = $ws.Cells.Item($r, 4).Value2   # D: 总天
    $remainingDays = $ws.Cells.Item($r, 5).Value2   # E: 剩余
    $startDate     = $ws.Cells.Item($r, 6).Value2   # F: 开始时间

    $newRemaining = $remainingDays - 1

    if ($newRemaining -le 0) {
        $newRemaining = $totalDays
        $ws.Cells.Item($r, 6).Value = $startDate + $totalDays
    }

    $ws.Cells.Item($r, 5).Value = $newRemaining
}
